# Natmi following Dr Hou advice
# Adds the "sCs" sending-cluster rows (8-10) to the Bmp8a-Tgfbr2 LR-pair sheet
# and refreshes rows 2-7 to reflect the 3-cluster (ECs/FAPs/sCs) recomputation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{}
$rowData[2] = @("ECs", "Bmp8a", "Tgfbr2", "ECs", 3, 1, 0.3786073333333334, 1.135822, 0.3713290366620658, 0.3713290366620658, 3, 1, 27.592778, 82.778334, 0.2684079248986126, 0.2684079248986126, 10.44682809783867, 94.02145288054801, 0.09966765618506591, 0.09966765618506591)
$rowData[3] = @("ECs", "Bmp8a", "Tgfbr2", "FAPs", 3, 1, 0.3786073333333334, 1.135822, 0.3713290366620658, 0.3713290366620658, 3, 1, 65.63594833333333, 196.907845, 0.6384717294830753, 0.6384717294830752, 24.85025136928778, 223.65226232359, 0.2370830922449134, 0.2370830922449134)
$rowData[4] = @("ECs", "Bmp8a", "Tgfbr2", "sCs", 3, 1, 0.3786073333333334, 1.135822, 0.3713290366620658, 0.3713290366620658, 3, 1, 9.572925333333334, 28.718776, 0.09312034561831214, 0.09312034561831213, 3.624379732652445, 32.61941759387201, 0.03457828823208647, 0.03457828823208646)
$rowData[5] = @("FAPs", "Bmp8a", "Tgfbr2", "ECs", 3, 1, 0.4265683333333333, 1.279705, 0.4183680408212104, 0.4183680408212104, 3, 1, 27.592778, 82.778334, 0.2684079248986126, 0.2684079248986126, 11.77020532349666, 105.93184791147, 0.1122932976807191, 0.1122932976807191)
$rowData[6] = @("FAPs", "Bmp8a", "Tgfbr2", "FAPs", 3, 1, 0.4265683333333333, 1.279705, 0.4183680408212104, 0.4183680408212104, 3, 1, 65.63594833333333, 196.907845, 0.6384717294830753, 0.6384717294830752, 27.99821708730277, 251.983953785725, 0.267116166583564, 0.267116166583564)
$rowData[7] = @("FAPs", "Bmp8a", "Tgfbr2", "sCs", 3, 1, 0.4265683333333333, 1.279705, 0.4183680408212104, 0.4183680408212104, 3, 1, 9.572925333333334, 28.718776, 0.09312034561831214, 0.09312034561831213, 4.083506804564444, 36.75156124108, 0.03895857655692724, 0.03895857655692723)
$rowData[8] = @("sCs", "Bmp8a", "Tgfbr2", "ECs", 2, 0.6666666666666666, 0.214425, 0.6432749999999999, 0.2103029225167239, 0.2103029225167238, 3, 1, 27.592778, 82.778334, 0.2684079248986126, 0.2684079248986126, 5.916581422649999, 53.24923280384999, 0.05644697103282756, 0.05644697103282755)
$rowData[9] = @("sCs", "Bmp8a", "Tgfbr2", "FAPs", 2, 0.6666666666666666, 0.214425, 0.6432749999999999, 0.2103029225167239, 0.2103029225167238, 3, 1, 65.63594833333333, 196.907845, 0.6384717294830753, 0.6384717294830752, 14.073988221375, 126.665893992375, 0.1342724706545979, 0.1342724706545978)
$rowData[10] = @("sCs", "Bmp8a", "Tgfbr2", "sCs", 2, 0.6666666666666666, 0.214425, 0.6432749999999999, 0.2103029225167239, 0.2103029225167238, 3, 1, 9.572925333333334, 28.718776, 0.09312034561831214, 0.09312034561831213, 2.0526745146, 18.4740706314, 0.01958348082929845, 0.01958348082929844)

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
